# C5-PowerPoint.pptx edit
#
# 1) Slide 6 table: switch the table's style from the custom "Table_0"
#    style to the built-in PowerPoint table style
#    {27F4BA68-CC3F-4165-9934-C4B8783AA56D}.
#
# 2) The deck's theme (ppt/theme/theme1.xml, the "Integral" theme used by
#    the Slide Master) is recoloured to the stock "Office" colour scheme
#    (the palette that used to live only in ppt/theme/theme2.xml, the
#    Notes Master's theme). Font scheme / format scheme are already
#    identical between the two themes, so only the 12 theme colours
#    change.

$p = $ppt.ActivePresentation

# --- 1. Table style -------------------------------------------------
$s6 = $p.Slides.Item(6)
$tableShape = $s6.Shapes.Item(2)
$tableShape.Table.ApplyStyle("{27F4BA68-CC3F-4165-9934-C4B8783AA56D}", $true)

# --- 2. Theme colours -------------------------------------------------
# Any slide exposes the deck's (single) ThemeColorScheme; slide 1 is as
# good as any other since they all share the one Slide Master / theme.
$tcs = $p.Slides.Item(1).ThemeColorScheme

$tcs.Colors(1).RGB  = 0         # dk1      000000
$tcs.Colors(2).RGB  = 16777215  # lt1      FFFFFF
$tcs.Colors(3).RGB  = 6968388   # dk2      44546A
$tcs.Colors(4).RGB  = 15132391  # lt2      E7E6E6
$tcs.Colors(5).RGB  = 13998939  # accent1  5B9BD5
$tcs.Colors(6).RGB  = 3243501   # accent2  ED7D31
$tcs.Colors(7).RGB  = 10855845  # accent3  A5A5A5
$tcs.Colors(8).RGB  = 49407     # accent4  FFC000
$tcs.Colors(9).RGB  = 12874308  # accent5  4472C4
$tcs.Colors(10).RGB = 4697456   # accent6  70AD47
$tcs.Colors(11).RGB = 12673797  # hlink    0563C1
$tcs.Colors(12).RGB = 7491477   # folHlink 954F72
